$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.070014146749288
$ws.Range("D2").Value = 1.06976653682712
$ws.Range("E2").Value = 1.073882399047981
$ws.Range("F2").Value = 1.082995814583443
$ws.Range("I2").Value = 1.051784793677603
$ws.Range("J2").Value = 1.074945950641331
$ws.Range("K2").Value = 1.072468166773134
$ws.Range("L2").Value = 1.07657308269903
$ws.Range("M2").Value = 1.085662580470776
$ws.Range("N2").Value = 1.076472497189266
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.071523160222331
$ws.Range("D3").Value = 1.070920030451972
$ws.Range("E3").Value = 1.07519063302456
$ws.Range("F3").Value = 1.084296318605072
$ws.Range("I3").Value = 1.052191866972033
$ws.Range("J3").Value = 1.076109709201634
$ws.Range("K3").Value = 1.073437057980772
$ws.Range("L3").Value = 1.077697129149831
$ws.Range("M3").Value = 1.086780657664182
$ws.Range("N3").Value = 1.077637908420209
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.072498501957363
$ws.Range("D4").Value = 1.071665239141248
$ws.Range("E4").Value = 1.076036295216696
$ws.Range("F4").Value = 1.085136950566755
$ws.Range("I4").Value = 1.052453160705257
$ws.Range("J4").Value = 1.076861207369685
$ws.Range("K4").Value = 1.074062234482782
$ws.Range("L4").Value = 1.078423057379435
$ws.Range("M4").Value = 1.087502690436903
$ws.Range("N4").Value = 1.078390473801857
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.072908280971594
$ws.Range("D5").Value = 1.071978246726999
$ws.Range("E5").Value = 1.076391612266518
$ws.Range("F5").Value = 1.085490145871162
$ws.Range("I5").Value = 1.052562505358848
$ws.Range("J5").Value = 1.07717677531045
$ws.Range("K5").Value = 1.074324640962951
$ws.Range("L5").Value = 1.078727905683463
$ws.Range("M5").Value = 1.087805893002075
$ws.Range("N5").Value = 1.078706489885294
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.072977069926967
$ws.Range("D6").Value = 1.072030785884931
$ws.Range("E6").Value = 1.076451260012431
$ws.Range("F6").Value = 1.085549436947364
$ws.Range("I6").Value = 1.052580835352088
$ws.Range("J6").Value = 1.077229739472235
$ws.Range("K6").Value = 1.074368675781454
$ws.Range("L6").Value = 1.078779071719362
$ws.Range("M6").Value = 1.087856782230989
$ws.Range("N6").Value = 1.078759529262263
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.072503978437145
$ws.Range("D7").Value = 1.071669422651835
$ws.Range("E7").Value = 1.07604104375782
$ws.Range("F7").Value = 1.085141670784952
$ws.Range("I7").Value = 1.052454623748909
$ws.Range("J7").Value = 1.07686542542007
$ws.Range("K7").Value = 1.074065742407631
$ws.Range("L7").Value = 1.078427132076143
$ws.Range("M7").Value = 1.087506743175205
$ws.Range("N7").Value = 1.078394697842357
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.070524354719213
$ws.Range("D8").Value = 1.070156611933517
$ws.Range("E8").Value = 1.074324702289432
$ws.Range("F8").Value = 1.083435511199689
$ws.Range("I8").Value = 1.051922804011985
$ws.Range("J8").Value = 1.07533956825027
$ws.Range("K8").Value = 1.072795975016875
$ws.Range("L8").Value = 1.076953253046286
$ws.Range("M8").Value = 1.08604074018111
$ws.Range("N8").Value = 1.076866673780358
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.067027357645081
$ws.Range("D9").Value = 1.067481634874247
$ws.Range("E9").Value = 1.071293543358975
$ws.Range("F9").Value = 1.080422086462373
$ws.Range("I9").Value = 1.050969425320523
$ws.Range("J9").Value = 1.072638870865947
$ws.Range("K9").Value = 1.070544823408709
$ws.Range("L9").Value = 1.07434513008538
$ws.Range("M9").Value = 1.083446252345462
$ws.Range("N9").Value = 1.074162141095975
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.064689801945728
$ws.Range("D10").Value = 1.065691865422584
$ws.Range("E10").Value = 1.069267931590954
$ws.Range("F10").Value = 1.078408171883951
$ws.Range("I10").Value = 1.050322797385038
$ws.Range("J10").Value = 1.070830076717234
$ws.Range("K10").Value = 1.069034628117904
$ws.Range("L10").Value = 1.072598732068127
$ws.Range("M10").Value = 1.081708792899329
$ws.Range("N10").Value = 1.0723507782521
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.063676046321728
$ws.Range("D11").Value = 1.064915287087681
$ws.Range("E11").Value = 1.06838960335434
$ws.Range("F11").Value = 1.077534881894925
$ws.Range("I11").Value = 1.050040154047087
$ws.Range("J11").Value = 1.070044805910135
$ws.Range("K11").Value = 1.068378405872905
$ws.Range("L11").Value = 1.071840645063449
$ws.Range("M11").Value = 1.080954543090285
$ws.Range("N11").Value = 1.071564392270383
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.063299247187841
$ws.Range("D12").Value = 1.064626586373036
$ws.Range("E12").Value = 1.068063163140597
$ws.Range("F12").Value = 1.077210309351166
$ws.Range("I12").Value = 1.049934767296726
$ws.Range("J12").Value = 1.069752807333647
$ws.Range("K12").Value = 1.068134305683095
$ws.Range("L12").Value = 1.07155876905765
$ws.Range("M12").Value = 1.080674087612174
$ws.Range("N12").Value = 1.071271979022423
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.063380083103936
$ws.Range("D13").Value = 1.064688524823033
$ws.Range("E13").Value = 1.068133194368573
$ws.Range("F13").Value = 1.077279940123644
$ws.Range("I13").Value = 1.049957391271801
$ws.Range("J13").Value = 1.069815456325313
$ws.Range("K13").Value = 1.068186681919831
$ws.Range("L13").Value = 1.071619245572759
$ws.Range("M13").Value = 1.080734259631942
$ws.Range("N13").Value = 1.071334716982842
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.063644905035343
$ws.Range("D14").Value = 1.064891428023025
$ws.Range("E14").Value = 1.068362623614986
$ws.Range("F14").Value = 1.07750805661298
$ws.Range("I14").Value = 1.050031450928151
$ws.Range("J14").Value = 1.070020675657163
$ws.Range("K14").Value = 1.068358235643595
$ws.Range("L14").Value = 1.071817351027104
$ws.Range("M14").Value = 1.080931366571394
$ws.Range("N14").Value = 1.071540227749683
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.063808037841725
$ws.Range("D15").Value = 1.065016410770218
$ws.Range("E15").Value = 1.068503957190091
$ws.Range("F15").Value = 1.07764858086454
$ws.Range("I15").Value = 1.050077028393087
$ws.Range("J15").Value = 1.070147076270961
$ws.Range("K15").Value = 1.06846388902756
$ws.Range("L15").Value = 1.071939371909487
$ws.Range("M15").Value = 1.081052771609285
$ws.Range("N15").Value = 1.071666807866851
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.064757047627759
$ws.Range("D16").Value = 1.065743370300484
$ws.Range("E16").Value = 1.069326197026683
$ws.Range("F16").Value = 1.078466102406645
$ws.Range("I16").Value = 1.050341499471207
$ws.Range("J16").Value = 1.070882148839985
$ws.Range("K16").Value = 1.0690781306059
$ws.Range("L16").Value = 1.072649003618119
$ws.Range("M16").Value = 1.081758809105228
$ws.Range("N16").Value = 1.072402924323237
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.065351908135489
$ws.Range("D17").Value = 1.066198941823191
$ws.Range("E17").Value = 1.069841634196852
$ws.Range("F17").Value = 1.078978572942866
$ws.Range("I17").Value = 1.050506684228331
$ws.Range("J17").Value = 1.071342687181122
$ws.Range("K17").Value = 1.069462809560941
$ws.Range("L17").Value = 1.073093628251578
$ws.Range("M17").Value = 1.082201170433633
$ws.Range("N17").Value = 1.072864116681643
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.06569872843854
$ws.Range("D18").Value = 1.066464515307926
$ws.Range("E18").Value = 1.070142162190979
$ws.Range("F18").Value = 1.079277367955157
$ws.Range("I18").Value = 1.050602778260494
$ws.Range("J18").Value = 1.071611113912563
$ws.Range("K18").Value = 1.069686964919002
$ws.Range("L18").Value = 1.0733527887646
$ws.Range("M18").Value = 1.082459007574784
$ws.Range("N18").Value = 1.073132924609832
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.065816959697688
$ws.Range("D19").Value = 1.066555043090603
$ws.Range("E19").Value = 1.070244614671332
$ws.Range("F19").Value = 1.079379229085722
$ws.Range("I19").Value = 1.050635500587562
$ws.Range("J19").Value = 1.07170260717067
$ws.Range("K19").Value = 1.06976335870062
$ws.Range("L19").Value = 1.073441125180494
$ws.Range("M19").Value = 1.082546892201484
$ws.Range("N19").Value = 1.073224547798859
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.065288100988677
$ws.Range("D20").Value = 1.066150079242926
$ws.Range("E20").Value = 1.069786344893068
$ws.Range("F20").Value = 1.078923602197348
$ws.Range("I20").Value = 1.050488987906335
$ws.Range("J20").Value = 1.07129329624585
$ws.Range("K20").Value = 1.069421560081808
$ws.Range("L20").Value = 1.073045943051126
$ws.Range("M20").Value = 1.082153728396309
$ws.Range("N20").Value = 1.072814655605578
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.063566928409065
$ws.Range("D21").Value = 1.06483168489074
$ws.Range("E21").Value = 1.06829506771327
$ws.Range("F21").Value = 1.077440887387248
$ws.Range("I21").Value = 1.050009653278613
$ws.Range("J21").Value = 1.069960252395321
$ws.Range("K21").Value = 1.068307727049259
$ws.Range("L21").Value = 1.071759021923542
$ws.Range("M21").Value = 1.080873331635702
$ws.Range("N21").Value = 1.071479718679879
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.062483336491475
$ws.Range("D22").Value = 1.064001338305411
$ws.Range("E22").Value = 1.067356340303403
$ws.Range("F22").Value = 1.076507521853168
$ws.Range("I22").Value = 1.049705958595898
$ws.Range("J22").Value = 1.069120294961959
$ws.Range("K22").Value = 1.067605388434837
$ws.Range("L22").Value = 1.070948210221033
$ws.Range("M22").Value = 1.080066594528355
$ws.Range("N22").Value = 1.070638568410617
$ws.Range("B23").Value = 1.019999999999999
$ws.Range("C23").Value = 1.063057906578087
$ws.Range("D23").Value = 1.064441657030321
$ws.Range("E23").Value = 1.067854083945925
$ws.Range("F23").Value = 1.077002425029135
$ws.Range("I23").Value = 1.049867173423835
$ws.Range("J23").Value = 1.069565746769
$ws.Range("K23").Value = 1.067977905233056
$ws.Range("L23").Value = 1.071378197307442
$ws.Range("M23").Value = 1.080494424131306
$ws.Range("N23").Value = 1.071084652810318
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.065316933161029
$ws.Range("D24").Value = 1.066172158614858
$ws.Range("E24").Value = 1.069811328113834
$ws.Range("F24").Value = 1.078948441483379
$ws.Range("I24").Value = 1.050496984901956
$ws.Range("J24").Value = 1.071315614494012
$ws.Range("K24").Value = 1.069440199630144
$ws.Range("L24").Value = 1.073067490501775
$ws.Range("M24").Value = 1.082175165981985
$ws.Range("N24").Value = 1.072837005548212
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.067932482052708
$ws.Range("D25").Value = 1.068174298619768
$ws.Range("E25").Value = 1.072078001161453
$ws.Range("F25").Value = 1.081201982710330
$ws.Range("I25").Value = 1.05121783393509
$ws.Range("J25").Value = 1.073338511843913
$ws.Range("K25").Value = 1.071128443746778
$ws.Range("L25").Value = 1.075020719960026
$ws.Range("M25").Value = 1.084118344807403
$ws.Range("N25").Value = 1.074862775644379
